$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.570.28'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '1.662.52'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("D4").Formula = '="0.9995"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Formula = '="235.64"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Formula = '="0.4792"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Formula = '="0.2616"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Formula = '="0.06146"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +1.88%  '
$ws.Range("D10").Formula = '="0.07078"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("D11").Value = '1.665.95'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Formula = '="14.73"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +1.89%  '
$ws.Range("D13").Formula = '="0.5902"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -4.48%  '
$ws.Range("D14").Formula = '="4.373"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -4.19%  '
$ws.Range("D15").Formula = '="74.35"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Formula = '="1.000"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").Value = '25.562.49'
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("D19").Formula = '="0.000006750"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +2.85%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '1.881.10'
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Formula = '="8.645"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +2.06%  '
$ws.Range("D24").Formula = '="5.287"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").Formula = '="1.401"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").Formula = '="104.61"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +2.78%  '
$ws.Range("D29").Formula = '="1.683"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").Formula = '="3.956"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +4.53%  '
$ws.Range("D31").Formula = '="3.653"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +2.91%  '
$ws.Range("D32").Formula = '="0.07611"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -3.89%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").Formula = '="0.04311"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -5.09%  '
$ws.Range("D35").Formula = '="2.617"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("E36").Value = '  +5.69%  '
$ws.Range("D37").Formula = '="0.9485"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("D38").Formula = '="2.605"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("D39").Formula = '="0.8520"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +1.72%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -2.62%  '
$ws.Range("D42").Formula = '="1.874"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +2.69%  '
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("D44").Formula = '="0.3755"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +1.28%  '
$ws.Range("D45").Formula = '="4.681"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("D46").Formula = '="0.1119"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Formula = '="6.204"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +2.73%  '
$ws.Range("D48").Formula = '="0.05259"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +2.08%  '
$ws.Range("D49").Formula = '="29.44"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Formula = '="7.348"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +1.05%  '

$excel.CutCopyMode = 0
